$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: update B7 value and add C7/D7 text, mirroring the style pattern of row 4
$ws.Range("B7").Value = 50
$ws.Range("D7").Value = "Email.qualquer2.com"
$ws.Range("C7").Value = "Nome"

# Copy formatting from row 4's C/D cells so the new cells match the workbook's existing style
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D4").Copy()
$ws.Range("D7").PasteSpecial(-4122) # xlPasteFormats

# Update the active selection to B7 to match the saved view state
$ws.Range("B7").Select()
